# Apply updated two-digit division answers to the answer table.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$mismatches = 0

$cell = $t.Cell(1, 1)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "46÷7=6, 4") {
    Write-Host "Mismatch at Cell(1,1): expected '46÷7=6, 4' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "47÷5=9, 2"
}

$cell = $t.Cell(1, 2)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "36÷2=18, 0") {
    Write-Host "Mismatch at Cell(1,2): expected '36÷2=18, 0' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "90÷6=15, 0"
}

$cell = $t.Cell(1, 3)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "54÷8=6, 6") {
    Write-Host "Mismatch at Cell(1,3): expected '54÷8=6, 6' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "45÷2=22, 1"
}

$cell = $t.Cell(1, 4)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "71÷6=11, 5") {
    Write-Host "Mismatch at Cell(1,4): expected '71÷6=11, 5' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "21÷9=2, 3"
}

$cell = $t.Cell(1, 5)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "15÷5=3, 0") {
    Write-Host "Mismatch at Cell(1,5): expected '15÷5=3, 0' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "65÷9=7, 2"
}

$cell = $t.Cell(5, 1)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "73÷8=9, 1") {
    Write-Host "Mismatch at Cell(5,1): expected '73÷8=9, 1' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "23÷3=7, 2"
}

$cell = $t.Cell(5, 2)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "94÷5=18, 4") {
    Write-Host "Mismatch at Cell(5,2): expected '94÷5=18, 4' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "25÷4=6, 1"
}

$cell = $t.Cell(5, 3)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "38÷6=6, 2") {
    Write-Host "Mismatch at Cell(5,3): expected '38÷6=6, 2' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "82÷4=20, 2"
}

$cell = $t.Cell(5, 4)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "66÷7=9, 3") {
    Write-Host "Mismatch at Cell(5,4): expected '66÷7=9, 3' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "61÷9=6, 7"
}

$cell = $t.Cell(5, 5)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "96÷4=24, 0") {
    Write-Host "Mismatch at Cell(5,5): expected '96÷4=24, 0' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "30÷3=10, 0"
}

$cell = $t.Cell(9, 1)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "50÷4=12, 2") {
    Write-Host "Mismatch at Cell(9,1): expected '50÷4=12, 2' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "69÷2=34, 1"
}

$cell = $t.Cell(9, 2)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "50÷5=10, 0") {
    Write-Host "Mismatch at Cell(9,2): expected '50÷5=10, 0' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "11÷2=5, 1"
}

$cell = $t.Cell(9, 3)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "23÷2=11, 1") {
    Write-Host "Mismatch at Cell(9,3): expected '23÷2=11, 1' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "46÷8=5, 6"
}

$cell = $t.Cell(9, 4)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "77÷9=8, 5") {
    Write-Host "Mismatch at Cell(9,4): expected '77÷9=8, 5' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "58÷2=29, 0"
}

$cell = $t.Cell(9, 5)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "23÷2=11, 1") {
    Write-Host "Mismatch at Cell(9,5): expected '23÷2=11, 1' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "99÷3=33, 0"
}

$cell = $t.Cell(13, 1)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "69÷8=8, 5") {
    Write-Host "Mismatch at Cell(13,1): expected '69÷8=8, 5' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "55÷2=27, 1"
}

$cell = $t.Cell(13, 2)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "57÷8=7, 1") {
    Write-Host "Mismatch at Cell(13,2): expected '57÷8=7, 1' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "53÷6=8, 5"
}

$cell = $t.Cell(13, 3)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "26÷7=3, 5") {
    Write-Host "Mismatch at Cell(13,3): expected '26÷7=3, 5' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "94÷2=47, 0"
}

$cell = $t.Cell(13, 4)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "82÷9=9, 1") {
    Write-Host "Mismatch at Cell(13,4): expected '82÷9=9, 1' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "36÷4=9, 0"
}

$cell = $t.Cell(13, 5)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "21÷6=3, 3") {
    Write-Host "Mismatch at Cell(13,5): expected '21÷6=3, 3' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "95÷3=31, 2"
}

$cell = $t.Cell(17, 1)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "53÷7=7, 4") {
    Write-Host "Mismatch at Cell(17,1): expected '53÷7=7, 4' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "59÷2=29, 1"
}

$cell = $t.Cell(17, 2)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "76÷7=10, 6") {
    Write-Host "Mismatch at Cell(17,2): expected '76÷7=10, 6' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "58÷7=8, 2"
}

$cell = $t.Cell(17, 3)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "96÷8=12, 0") {
    Write-Host "Mismatch at Cell(17,3): expected '96÷8=12, 0' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "69÷5=13, 4"
}

$cell = $t.Cell(17, 4)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "42÷3=14, 0") {
    Write-Host "Mismatch at Cell(17,4): expected '42÷3=14, 0' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "54÷6=9, 0"
}

$cell = $t.Cell(17, 5)
$cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)
if ($cellText -ne "74÷7=10, 4") {
    Write-Host "Mismatch at Cell(17,5): expected '74÷7=10, 4' but found '$cellText'"
    $mismatches++
} else {
    $cell.Range.Text = "59÷9=6, 5"
}

Write-Host "Done. Mismatches: $mismatches"